$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 522.53845
$ws.Range("I28").Value = 305
$ws.Range("J28").Value = 709
$ws.Range("K28").Value = 305
$ws.Range("L28").Value = 709
$ws.Range("M28").Value = 180
$ws.Range("N28").Value = -1679
$ws.Range("H98").Value = 12700548
$ws.Range("I98").Value = 16463393
$ws.Range("J98").Value = 945.25
$ws.Range("K98").Value = 16463393
$ws.Range("L98").Value = 945.25
$ws.Range("M98").Value = -16461895
$ws.Range("N98").Value = -3941.25
$ws.Range("H122").Value = 12700548
$ws.Range("I122").Value = 16463393
$ws.Range("J122").Value = 945.25
$ws.Range("K122").Value = 49390179
$ws.Range("L122").Value = 2835.75
$ws.Range("M122").Value = -49387729
$ws.Range("N122").Value = -7735.75
$ws.Range("H123").Value = 39987.5
$ws.Range("J123").Value = 39987.5
$ws.Range("L123").Value = 39987.5
$ws.Range("N123").Value = -49787.5
$ws.Range("H137").Value = 593017.75
$ws.Range("I137").Value = 796163.5600000001
$ws.Range("J137").Value = 85153.336
$ws.Range("K137").Value = 2388490.68
$ws.Range("L137").Value = 255460.008
$ws.Range("M137").Value = -2385940.68
$ws.Range("N137").Value = -260560.008
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5878.409
$ws.Range("I61").Value = 7462.3335
$ws.Range("J61").Value = 2484.2856
$ws.Range("K61").Value = 7462.3335
$ws.Range("L61").Value = 2484.2856
$ws.Range("M61").Value = -7250.3335
$ws.Range("N61").Value = -2908.2856
$ws.Range("H74").Value = 3741.081
$ws.Range("I74").Value = 742.6818
$ws.Range("J74").Value = 8138.7334
$ws.Range("K74").Value = 742.6818
$ws.Range("L74").Value = 8138.7334
$ws.Range("M74").Value = 131.3182
$ws.Range("N74").Value = -9886.733400000001
$ws.Range("H77").Value = 3741.081
$ws.Range("I77").Value = 742.6818
$ws.Range("J77").Value = 8138.7334
$ws.Range("K77").Value = 3713.409
$ws.Range("L77").Value = 40693.667
$ws.Range("M77").Value = 654.5910000000003
$ws.Range("N77").Value = -49429.667
$ws.Range("H132").Value = 5954557.5
$ws.Range("I132").Value = 10417949
$ws.Range("J132").Value = 3367.889
$ws.Range("K132").Value = 31253847
$ws.Range("L132").Value = 10103.667
$ws.Range("M132").Value = -31251317
$ws.Range("N132").Value = -15163.667
$ws.Range("H136").Value = 5878.409
$ws.Range("I136").Value = 7462.3335
$ws.Range("J136").Value = 2484.2856
$ws.Range("K136").Value = 22387.0005
$ws.Range("L136").Value = 7452.8568
$ws.Range("M136").Value = -19837.0005
$ws.Range("N136").Value = -12552.8568
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 27780164
$ws.Range("I134").Value = 47621280
$ws.Range("J134").Value = 2600
$ws.Range("K134").Value = 142863840
$ws.Range("L134").Value = 7800
$ws.Range("M134").Value = -142861305
$ws.Range("N134").Value = -12870
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10000.869
$ws.Range("I31").Value = 1065.0588
$ws.Range("J31").Value = 35319
$ws.Range("K31").Value = 1065.0588
$ws.Range("L31").Value = 35319
$ws.Range("M31").Value = -770.0588
$ws.Range("N31").Value = -35909
$ws.Range("H34").Value = 10000.869
$ws.Range("I34").Value = 1065.0588
$ws.Range("J34").Value = 35319
$ws.Range("K34").Value = 1065.0588
$ws.Range("L34").Value = 35319
$ws.Range("M34").Value = -863.0588
$ws.Range("N34").Value = -35723
$ws.Range("H58").Value = 4647938
$ws.Range("I58").Value = 7193648
$ws.Range("K58").Value = 7193648
$ws.Range("M58").Value = -7193445
$ws.Range("H122").Value = 3664476
$ws.Range("I122").Value = 6212200
$ws.Range("J122").Value = 2122.5625
$ws.Range("K122").Value = 18636600
$ws.Range("L122").Value = 6367.6875
$ws.Range("M122").Value = -18634150
$ws.Range("N122").Value = -11267.6875
$ws.Range("H132").Value = 14499402
$ws.Range("I132").Value = 83334290
$ws.Range("J132").Value = 7847.684
$ws.Range("K132").Value = 250002870
$ws.Range("L132").Value = 23543.052
$ws.Range("M132").Value = -250000340
$ws.Range("N132").Value = -28603.052
$ws.Range("H134").Value = 14882724
$ws.Range("I134").Value = 16668227
$ws.Range("J134").Value = 10418967
$ws.Range("K134").Value = 50004681
$ws.Range("L134").Value = 31256901
$ws.Range("M134").Value = -50002146
$ws.Range("N134").Value = -31261971
$ws.Range("H136").Value = 4647938
$ws.Range("I136").Value = 7193648
$ws.Range("K136").Value = 21580944
$ws.Range("M136").Value = -21578394
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1733.4
$ws.Range("I113").Value = 1523.1538
$ws.Range("J113").Value = 1961.1666
$ws.Range("K113").Value = 1523.1538
$ws.Range("L113").Value = 1961.1666
$ws.Range("M113").Value = 646.8462
$ws.Range("N113").Value = -6301.1666
$ws.Range("H122").Value = 37038188
$ws.Range("I122").Value = 55556460
$ws.Range("J122").Value = 1641.3334
$ws.Range("K122").Value = 166669380
$ws.Range("L122").Value = 4924.0002
$ws.Range("M122").Value = -166666930
$ws.Range("N122").Value = -9824.0002
$ws.Range("H132").Value = 55557372
$ws.Range("I132").Value = 76924056
$ws.Range("J132").Value = 3999.6
$ws.Range("K132").Value = 230772168
$ws.Range("L132").Value = 11998.8
$ws.Range("M132").Value = -230769638
$ws.Range("N132").Value = -17058.8
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1731.25
$ws.Range("I40").Value = 1681.6666
$ws.Range("J40").Value = 1780.8334
$ws.Range("K40").Value = 1681.6666
$ws.Range("L40").Value = 1780.8334
$ws.Range("M40").Value = -1545.6666
$ws.Range("N40").Value = -2052.8334
$ws.Range("H122").Value = 48079644
$ws.Range("I122").Value = 71430560
$ws.Range("J122").Value = 20836910
$ws.Range("K122").Value = 214291680
$ws.Range("L122").Value = 62510730
$ws.Range("M122").Value = -214289230
$ws.Range("N122").Value = -62515630
$ws.Range("H132").Value = 10002344
$ws.Range("I132").Value = 18183082
$ws.Range("J132").Value = 3666.4443
$ws.Range("K132").Value = 54549246
$ws.Range("L132").Value = 10999.3329
$ws.Range("M132").Value = -54546716
$ws.Range("N132").Value = -16059.3329
$ws.Range("H136").Value = 2413.647
$ws.Range("I136").Value = 2313.2964
$ws.Range("J136").Value = 2800.7144
$ws.Range("K136").Value = 6939.889200000001
$ws.Range("L136").Value = 8402.143199999999
$ws.Range("M136").Value = -4389.889200000001
$ws.Range("N136").Value = -13502.1432
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4260.9644
$ws.Range("I122").Value = 5605.35
$ws.Range("K122").Value = 16816.05
$ws.Range("M122").Value = -14366.05
$ws.Range("H132").Value = 95226260
$ws.Range("I132").Value = 133334504
$ws.Range("K132").Value = 400003512
$ws.Range("M132").Value = -400000982
$ws.Range("H136").Value = 19910322
$ws.Range("I136").Value = 12718051
$ws.Range("J136").Value = 33335894
$ws.Range("K136").Value = 38154153
$ws.Range("L136").Value = 100007682
$ws.Range("M136").Value = -38151603
$ws.Range("N136").Value = -100012782

Write-Host "Done updating cells"
